# Capítulo 11. Presupuesto final.xlsx
# "Teminados todos los puntos importantes de la documentacion"
#
# The underlying data edit is a single cell: on "Planificacion_Final",
# D7 (hours for the "Documentación" task) goes from 73 to 82. Every other
# changed cell in the diff (D2 on that sheet, and B3/F3/I3/M3/N3/M4/N4/M6/N6/
# B8/F8 on "Presupuesto final", plus the chart cache) is a downstream
# formula recalculation that Excel performs automatically once the input
# cell changes - no manual edits are required for those.
#
# The diff also shows view/selection-state changes: "Presupuesto final"
# becomes the active tab with L3:N6 selected (active cell N6), while
# "Planificacion_Final" is no longer the active tab and has D2:D7 selected.

$wb = $excel.ActiveWorkbook

# --- The actual data edit ---
$wsPlan = $wb.Worksheets.Item("Planificacion_Final")
$wsPlan.Range("D7").Value = 82

# --- View state: Planificacion_Final selection becomes D2:D7 ---
$wsPlan.Select()
$wsPlan.Range("D2:D7").Select()

# --- View state: Presupuesto final becomes the active tab, selection L3:N6 ---
$wsPres = $wb.Worksheets.Item("Presupuesto final")
$wsPres.Select()
$wsPres.Range("L3:N6").Select()
